# edit.ps1 - applies the diff to ToDo Documentation.docx via Word COM interop
$d = $word.ActiveDocument

# --- 1. Paragraph 45 (last paragraph, was an empty underline-formatted paragraph
#        right after "Test Case") -> replaced by the large "Test Case" write-up
#        block (19 paragraphs) ending with the relocated _GoBack bookmark.
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIndex)
$pLast.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Name of task:</w:t></w:r><w:r><w:t xml:space="preserve"> Learning how to use XAMPP and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>myPHPAdmin</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Estimated No of Hours</w:t></w:r><w:r><w:t>: 1</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Actual time spent:</w:t></w:r><w:r><w:t xml:space="preserve"> 2</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Notes:</w:t></w:r><w:r><w:t xml:space="preserve"> Had to learn how to use both programs, and add the data to the databases through the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and MySQL coding</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Issues:</w:t></w:r><w:r><w:t xml:space="preserve"> Finding the proper material to learn how these programs function</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Name of task:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Learning PHP and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mySql</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Estimated No of Hours:</w:t></w:r><w:r><w:t xml:space="preserve"> 3</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Actual time spent:</w:t></w:r><w:r><w:t xml:space="preserve"> 6</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Notes:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Using an HTML interface for the user and a PHP and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mySql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> backend, I was able to connect to the database, allowing me to store and view information that was used in my code.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Issues:</w:t></w:r><w:r><w:t xml:space="preserve"> My biggest issue is being able to use the add and remove button to use for the tasks. I am able to hardcode the tasks into the database, but struggled with being able to use the interface.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Name of task</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> Completing documentation</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Estimated No of Hours:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>1</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Actual time spent:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>2</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Notes: </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Created detailed information and diagrams about how my </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>todo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> list application was setup.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Issues:</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Making sure the diagrams properly portrayed what was needed seemed to be the most difficult for me.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# --- 2. Paragraph 44 ("Test Case" heading) -> add <w:lastRenderedPageBreak/>
#        immediately before the run text, keep pPr/rPr untouched otherwise.
$p44 = $d.Paragraphs(44)
$p44.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B32FA2" w:rsidRDefault="00192737" w:rsidP="00B32FA2"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r w:rsidRPr="00192737"><w:rPr><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Test Case</w:t></w:r></w:p>')

# --- 3. Paragraph 41 ("storeData" bullet) -> add <w:proofErr w:type="spellStart"/>
#        right before the "storeData" run.
$p41 = $d.Paragraphs(41)
$p41.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B32FA2" w:rsidRDefault="00B32FA2" w:rsidP="00B32FA2"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>storeData</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class, that would have the code pushed to the database and stored.</w:t></w:r></w:p>')

# --- 4. Paragraph 40 ("No classes were added...") -> merge the two runs into a
#        single run and drop the now-unneeded <w:proofErr w:type="spellStart"/>.
$p40 = $d.Paragraphs(40)
$p40.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B32FA2" w:rsidRDefault="00B32FA2" w:rsidP="000E1665"><w:r><w:t>No classes were added in this assignment, but if there was to be 3 classes added, I would have:</w:t></w:r></w:p>')

# --- 5. Paragraph 38 (empty paragraph right after "Database Design") -> filled
#        in with the new database-design description paragraph.
$p38 = $d.Paragraphs(38)
$p38.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B32FA2" w:rsidRDefault="00B32FA2" w:rsidP="000E1665"><w:r><w:t>My database design consists of 3 tables: tasks, completed and users. The tasks table stores the information that is to be entered by the user, the completed table will hold the information that the task has been completed, and the users table will hold the information pertaining to the user and their list.</w:t></w:r></w:p>')

# --- 6. Paragraph 37 ("Database Design" heading) -> give the paragraph mark
#        itself the single-underline formatting (pPr/rPr), matching the
#        pattern used by the other section headings.
$p37 = $d.Paragraphs(37)
$p37.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00B32FA2" w:rsidRDefault="00B32FA2" w:rsidP="000E1665"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Database Design</w:t></w:r></w:p>')

# --- 7. Paragraph 33 (the paragraph that used to hold the _GoBack bookmark)
#        -> bookmark removed from here (it is re-added at the very end of the
#        document in step 1 above), paragraph becomes empty.
$p33 = $d.Paragraphs(33)
$p33.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E10B78" w:rsidRDefault="00E10B78" w:rsidP="000E1665"/>')

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
